# Update "paises.xlsx" with the new country data scrape (1 Oct 2020, 22:54)
# - updates stats for several countries
# - Costa Rica overtakes Portugal and Etiopia, so those 3 rows are re-ranked
# - Islas Malvinas overtakes Montserrat, so those 2 rows swap
# - timestamp footer row updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (row 1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Octubre de 2020 a las 22:54"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 7484719
$ws.Range("C4").Value = 37437
$ws.Range("D4").Value = 4725243
$ws.Range("E4").Value = 2547081
$ws.Range("G4").Value = 655
$ws.Range("H4").Value = 212395

# Sudafrica (row 13)
$ws.Range("B13").Value = 676084
$ws.Range("C13").Value = 1745
$ws.Range("D13").Value = 609584
$ws.Range("E13").Value = 49634
$ws.Range("G13").Value = 132
$ws.Range("H13").Value = 16866

# Alemania (row 25)
$ws.Range("B25").Value = 295530
$ws.Range("C25").Value = 2619
$ws.Range("E25").Value = 28044
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 9586

# Israel (row 27)
$ws.Range("B27").Value = 253490
$ws.Range("C27").Value = 7996
$ws.Range("D27").Value = 179468
$ws.Range("E27").Value = 72400
$ws.Range("G27").Value = 53
$ws.Range("H27").Value = 1622

# Canada (row 29)
$ws.Range("B29").Value = 160280
$ws.Range("C29").Value = 1522
$ws.Range("D29").Value = 136089
$ws.Range("E29").Value = 14875

# Costa Rica now ranks above Portugal and Etiopia (rows 52-54 re-ordered)
$ws.Range("A52").Value = "Costa Rica"
$ws.Range("B52").Value = 76828
$ws.Range("C52").Value = 1068
$ws.Range("D52").Value = 39843
$ws.Range("E52").Value = 36068
$ws.Range("G52").Value = 13
$ws.Range("H52").Value = 917

$ws.Range("A53").Value = "Portugal"
$ws.Range("B53").Value = 76396
$ws.Range("C53").Value = 854
$ws.Range("D53").Value = 48937
$ws.Range("E53").Value = 25482
$ws.Range("G53").Value = 6
$ws.Range("H53").Value = 1977

$ws.Range("A54").Value = "Etiopia"
$ws.Range("B54").Value = 76098
$ws.Range("C54").Value = 730
$ws.Range("D54").Value = 31430
$ws.Range("E54").Value = 43463
$ws.Range("G54").Value = 7
$ws.Range("H54").Value = 1205

# Costa de Marfil (row 85)
$ws.Range("B85").Value = 19755
$ws.Range("C85").Value = 31
$ws.Range("D85").Value = 19320
$ws.Range("E85").Value = 315

# Luxemburgo (row 111)
$ws.Range("B111").Value = 8595
$ws.Range("C111").Value = 86
$ws.Range("D111").Value = 7280
$ws.Range("E111").Value = 1190
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 125

# Curazao (row 181)
$ws.Range("B181").Value = 399
$ws.Range("C181").Value = 7
$ws.Range("D181").Value = 185
$ws.Range("E181").Value = 213

# Nueva Caledonia now ranks above Santa Lucia (rows 207-208 swapped; stats tied)
$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("A208").Value = "Santa Lucia"

# Islas Malvinas now ranks above Montserrat (rows 215-216 swapped)
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
